# Add strings to relevant spreadsheets
# Updates the "Update"/"Results for this location"/"Language" block (rows 61-65)
# and appends a new row 66 ("Getting results for your location…") on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 61: "Update" / new Chinese label "最新信息" ------------------------
$ws.Cells.Item(61, 1).Value = "Update"
$ws.Cells.Item(61, 2).Value = "最新信息"

# --- Row 62: ellipsis dropped from the English string, Chinese gets a
#     full (double-character) ellipsis -------------------------------------
$ws.Cells.Item(62, 1).Value = "Results for this location"
$ws.Cells.Item(62, 2).Value = "该地方结果……"

# --- Row 63: unchanged text, kept for completeness -------------------------
$ws.Cells.Item(63, 1).Value = "Language"
$ws.Cells.Item(63, 2).Value = "语言"

# --- Row 64: unchanged text -------------------------------------------------
$ws.Cells.Item(64, 1).Value = "Your address, your city"
$ws.Cells.Item(64, 2).Value = "您的地址，所在城市"

# --- Row 65: unchanged text -------------------------------------------------
$ws.Cells.Item(65, 1).Value = "What to Expect at This Location"
$ws.Cells.Item(65, 2).Value = "对该地方有何期待"

# --- Row 66: brand-new row --------------------------------------------------
$ws.Cells.Item(66, 1).Value = "Getting results for your location…"
$ws.Cells.Item(66, 2).Value = "在该地方获得成果……"

# --- Formatting -------------------------------------------------------------
# Column A (English) reverts to the plain default/"Normal" style for every
# row in the block; column B (Chinese) keeps using the dedicated font, whose
# family is switched from "Microsoft YaHei" to "Calibri".
for ($r = 61; $r -le 66; $r++) {
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Font.Size = 11
    $ws.Cells.Item($r, 2).Font.Name = "Calibri"
}

$ws.Range("A61:B66").RowHeight = 16
